$wb = $excel.ActiveWorkbook

# Set Sheet2 data: A1 = 1020, B1 = "sheet2"
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("A1").Value = 1020
$ws2.Range("B1").Value = "sheet2"
$ws2.Range("B1").Select()

# Make Sheet2 the active/selected tab
$ws2.Activate()
